$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update threshold value in C9 from 0.8 to 0.7
$ws.Range("C9").Value = 0.7

# Update the selected cell on the sheet
$ws.Range("C11").Select()

# Update the workbook window position
$excel.ActiveWindow.Left = 760
$excel.ActiveWindow.Top = 480
